$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H5").Value = 274.1111
$ws.Range("I5").Value = 264.625
$ws.Range("K5").Value = 264.625
$ws.Range("M5").Value = -149.625
$ws.Range("H18").Value = 6187.3335
$ws.Range("I18").Value = 892.7143
$ws.Range("J18").Value = 13599.8
$ws.Range("K18").Value = 892.7143
$ws.Range("L18").Value = 13599.8
$ws.Range("M18").Value = -608.7143
$ws.Range("N18").Value = -14167.8
$ws.Range("H31").Value = 1235.375
$ws.Range("I31").Value = 1235.375
$ws.Range("K31").Value = 3706.125
$ws.Range("M31").Value = -3476.125
$ws.Range("H33").Value = 188.27586
$ws.Range("I33").Value = 191.78947
$ws.Range("J33").Value = 181.6
$ws.Range("K33").Value = 191.78947
$ws.Range("L33").Value = 181.6
$ws.Range("M33").Value = 37.21053000000001
$ws.Range("N33").Value = -639.6
$ws.Range("H47").Value = 7850
$ws.Range("I47").Value = 7850
$ws.Range("K47").Value = 7850
$ws.Range("M47").Value = -6878
$ws.Range("H69").Value = 9416
$ws.Range("H72").Value = 9416
$ws.Range("H96").Value = 604.8333
$ws.Range("I96").Value = 668.875
$ws.Range("K96").Value = 2006.625
$ws.Range("M96").Value = -633.625
$ws.Range("H98").Value = 2046.92
$ws.Range("I98").Value = 1987.1177
$ws.Range("K98").Value = 1987.1177
$ws.Range("M98").Value = -489.1177
$ws.Range("H112").Value = 1555.4524
$ws.Range("J112").Value = 1590.475
$ws.Range("L112").Value = 4771.424999999999
$ws.Range("N112").Value = -6987.424999999999
$ws.Range("H122").Value = 2046.92
$ws.Range("I122").Value = 1987.1177
$ws.Range("K122").Value = 5961.3531
$ws.Range("M122").Value = -3511.3531
$ws.Range("H125").Value = 3327.9167
$ws.Range("J125").Value = 3245.8333
$ws.Range("L125").Value = 29212.4997
$ws.Range("N125").Value = -34132.4997
$ws.Range("H138").Value = 2840.5962
$ws.Range("J138").Value = 3168.9119
$ws.Range("L138").Value = 9506.735700000001
$ws.Range("N138").Value = -19786.7357

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H30").Value = 12
$ws.Range("J30").Value = 0
$ws.Range("L30").Value = 0
$ws.Range("N30").ClearContents()
$ws.Range("H32").Value = 2610.7104
$ws.Range("I32").Value = 2607.1082
$ws.Range("J32").Value = 2744
$ws.Range("K32").Value = 2607.1082
$ws.Range("L32").Value = 2744
$ws.Range("M32").Value = -2320.1082
$ws.Range("N32").Value = -3318
$ws.Range("H45").Value = 5141.1816
$ws.Range("I45").Value = 6269.25
$ws.Range("K45").Value = 6269.25
$ws.Range("M45").Value = -5892.25
$ws.Range("H61").Value = 1525431.9
$ws.Range("I61").Value = 2789400.2
$ws.Range("K61").Value = 2789400.2
$ws.Range("M61").Value = -2789188.2
$ws.Range("H74").Value = 4202.22
$ws.Range("I74").Value = 3871.9768
$ws.Range("J74").Value = 6230.857
$ws.Range("K74").Value = 3871.9768
$ws.Range("L74").Value = 6230.857
$ws.Range("M74").Value = -2997.9768
$ws.Range("N74").Value = -7978.857
$ws.Range("H77").Value = 4202.22
$ws.Range("I77").Value = 3871.9768
$ws.Range("J77").Value = 6230.857
$ws.Range("K77").Value = 19359.884
$ws.Range("L77").Value = 31154.285
$ws.Range("M77").Value = -14991.884
$ws.Range("N77").Value = -39890.285
$ws.Range("H136").Value = 1525431.9
$ws.Range("I136").Value = 2789400.2
$ws.Range("K136").Value = 8368200.600000001
$ws.Range("M136").Value = -8365650.600000001

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H7").Value = 875
$ws.Range("I7").Value = 550
$ws.Range("J7").Value = 1200
$ws.Range("K7").Value = 550
$ws.Range("L7").Value = 1200
$ws.Range("M7").Value = -437
$ws.Range("N7").Value = -1426
$ws.Range("H20").Value = 1794.3077
$ws.Range("I20").Value = 2027.1666
$ws.Range("K20").Value = 2027.1666
$ws.Range("M20").Value = -1780.1666
$ws.Range("H94").Value = 1507.0769
$ws.Range("I94").Value = 1459.7
$ws.Range("J94").Value = 1665
$ws.Range("K94").Value = 1459.7
$ws.Range("L94").Value = 1665
$ws.Range("M94").Value = -1008.7
$ws.Range("N94").Value = -2567

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H107").Value = 1491.8
$ws.Range("I107").Value = 1215.3846
$ws.Range("J107").Value = 2005.1428
$ws.Range("K107").Value = 1215.3846
$ws.Range("L107").Value = 2005.1428
$ws.Range("M107").Value = 704.6153999999999
$ws.Range("N107").Value = -5845.1428
$ws.Range("H125").Value = 30000
$ws.Range("J125").Value = 30000
$ws.Range("L125").Value = 30000
$ws.Range("N125").Value = -34920

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H60").Value = 1689.4546
$ws.Range("I60").Value = 1358.4
$ws.Range("K60").Value = 4075.2
$ws.Range("M60").Value = -3824.2
$ws.Range("H132").Value = 2707.7273
$ws.Range("I132").Value = 2378.2
$ws.Range("J132").Value = 2982.3333
$ws.Range("K132").Value = 21403.8
$ws.Range("L132").Value = 26840.9997
$ws.Range("M132").Value = -18873.8
$ws.Range("N132").Value = -31900.9997
$ws.Range("H140").Value = 3796.4
$ws.Range("I140").Value = 2602.5652
$ws.Range("K140").Value = 7807.6956
$ws.Range("M140").Value = -2627.6956

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H5").Value = 10
$ws.Range("I5").Value = 10
$ws.Range("K5").Value = 10
$ws.Range("M5").Value = 102
$ws.Range("H21").Value = 920499.9399999999
$ws.Range("I21").Value = 1434814.1
$ws.Range("J21").Value = 20450
$ws.Range("K21").Value = 1434814.1
$ws.Range("L21").Value = 20450
$ws.Range("M21").Value = -1434641.1
$ws.Range("N21").Value = -20796
$ws.Range("H30").Value = 920499.9399999999
$ws.Range("I30").Value = 1434814.1
$ws.Range("J30").Value = 20450
$ws.Range("K30").Value = 1434814.1
$ws.Range("L30").Value = 20450
$ws.Range("M30").Value = -1434709.1
$ws.Range("N30").Value = -20660
$ws.Range("H49").Value = 0
$ws.Range("J49").Value = 0
$ws.Range("L49").ClearContents()
$ws.Range("N49").ClearContents()
$ws.Range("H126").Value = 3261.5557
$ws.Range("I126").Value = 2530.9
$ws.Range("K126").Value = 7592.700000000001
$ws.Range("M126").Value = -5122.700000000001

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3204.8
$ws.Range("I7").Value = 3204
$ws.Range("J7").Value = 3205
$ws.Range("K7").Value = 3204
$ws.Range("L7").Value = 3205
$ws.Range("M7").Value = -3092
$ws.Range("N7").Value = -3429
$ws.Range("H42").Value = 21998
$ws.Range("I42").Value = 21998
$ws.Range("K42").Value = 21998
$ws.Range("M42").Value = -21435
$ws.Range("H49").Value = 21998
$ws.Range("I49").Value = 21998
$ws.Range("K49").Value = 21998
$ws.Range("M49").Value = -21851
$ws.Range("H61").Value = 3326.5
$ws.Range("I61").Value = 1984.2142
$ws.Range("K61").Value = 1984.2142
$ws.Range("M61").Value = -1782.2142
$ws.Range("H113").Value = 3326.5
$ws.Range("I113").Value = 1984.2142
$ws.Range("K113").Value = 1984.2142
$ws.Range("M113").Value = 185.7858000000001
$ws.Range("H122").Value = 3123.3408
$ws.Range("I122").Value = 2850.6365
$ws.Range("K122").Value = 8551.9095
$ws.Range("M122").Value = -6101.9095
$ws.Range("H126").Value = 3204.8
$ws.Range("I126").Value = 3204
$ws.Range("J126").Value = 3205
$ws.Range("K126").Value = 9612
$ws.Range("L126").Value = 9615
$ws.Range("M126").Value = -7142
$ws.Range("N126").Value = -14555
$ws.Range("H130").Value = 63009
$ws.Range("J130").Value = 63009
$ws.Range("L130").Value = 63009
$ws.Range("N130").Value = -73049

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H70").Value = 48828.75
$ws.Range("I70").Value = 45000
$ws.Range("K70").Value = 45000
$ws.Range("M70").Value = -44685
$ws.Range("H73").Value = 48828.75
$ws.Range("I73").Value = 45000
$ws.Range("K73").Value = 45000
$ws.Range("M73").Value = -43908
$ws.Range("H74").Value = 179204.86
$ws.Range("J74").Value = 216687.2
$ws.Range("L74").Value = 216687.2
$ws.Range("N74").Value = -218559.2
$ws.Range("H77").Value = 179204.86
$ws.Range("J77").Value = 216687.2
$ws.Range("L77").Value = 650061.6000000001
$ws.Range("N77").Value = -659421.6000000001
$ws.Range("H81").Value = 1130.0769
$ws.Range("I81").Value = 1082.4445
$ws.Range("K81").Value = 2164.889
$ws.Range("M81").Value = -1103.889
$ws.Range("H84").Value = 1130.0769
$ws.Range("I84").Value = 1082.4445
$ws.Range("K84").Value = 10824.445
$ws.Range("M84").Value = -5520.445
$ws.Range("H126").Value = 2755.9412
$ws.Range("I126").Value = 2269.6155
$ws.Range("J126").Value = 4336.5
$ws.Range("K126").Value = 6808.8465
$ws.Range("L126").Value = 13009.5
$ws.Range("M126").Value = -4338.8465
$ws.Range("N126").Value = -17949.5
$ws.Range("H132").Value = 10067331
$ws.Range("I132").Value = 12581852
$ws.Range("J132").Value = 9249.25
$ws.Range("K132").Value = 37745556
$ws.Range("L132").Value = 27747.75
$ws.Range("M132").Value = -37743026
$ws.Range("N132").Value = -32807.75
